$wb = $excel.ActiveWorkbook

# ALC row 7
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(7, 8).Value = 5125
$ws.Cells.Item(7, 9).Value = 250
$ws.Cells.Item(7, 10).Value = 10000
$ws.Cells.Item(7, 11).Value = 250
$ws.Cells.Item(7, 12).Value = 10000
$ws.Cells.Item(7, 13).Value = -138
$ws.Cells.Item(7, 14).Value = -10224

# ALC row 14
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(14, 8).Value = 5125
$ws.Cells.Item(14, 9).Value = 250
$ws.Cells.Item(14, 10).Value = 10000
$ws.Cells.Item(14, 11).Value = 250
$ws.Cells.Item(14, 12).Value = 10000
$ws.Cells.Item(14, 13).Value = -59
$ws.Cells.Item(14, 14).Value = -10382

# ALC row 32
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 4356817
$ws.Cells.Item(32, 9).Value = 789.7778
$ws.Cells.Item(32, 10).Value = 9957423
$ws.Cells.Item(32, 11).Value = 789.7778
$ws.Cells.Item(32, 12).Value = 9957423
$ws.Cells.Item(32, 13).Value = -463.7778
$ws.Cells.Item(32, 14).Value = -9958075

# ALC row 54
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(54, 8).Value = 6000
$ws.Cells.Item(54, 10).Value = 6000
$ws.Cells.Item(54, 12).Value = 6000
$ws.Cells.Item(54, 14).Value = -6972

# ALC row 69
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(69, 8).Value = 3000
$ws.Cells.Item(69, 10).Value = 3000
$ws.Cells.Item(69, 12).Value = 9000
$ws.Cells.Item(69, 14).Value = -10748

# ALC row 72
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(72, 8).Value = 3000
$ws.Cells.Item(72, 10).Value = 3000
$ws.Cells.Item(72, 12).Value = 27000
$ws.Cells.Item(72, 14).Value = -35736

# ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 3651.6052
$ws.Cells.Item(74, 9).Value = 3202.0476
$ws.Cells.Item(74, 10).Value = 4206.9414
$ws.Cells.Item(74, 11).Value = 3202.0476
$ws.Cells.Item(74, 12).Value = 4206.9414
$ws.Cells.Item(74, 13).Value = -2266.0476
$ws.Cells.Item(74, 14).Value = -6078.9414

# ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(77, 8).Value = 3651.6052
$ws.Cells.Item(77, 9).Value = 3202.0476
$ws.Cells.Item(77, 10).Value = 4206.9414
$ws.Cells.Item(77, 11).Value = 16010.238
$ws.Cells.Item(77, 12).Value = 21034.707
$ws.Cells.Item(77, 13).Value = -11330.238
$ws.Cells.Item(77, 14).Value = -30394.707

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 6684.933
$ws.Cells.Item(86, 10).Value = 7232.5835
$ws.Cells.Item(86, 12).Value = 7232.5835
$ws.Cells.Item(86, 14).Value = -9478.583500000001

# ALC row 88
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(88, 8).Value = 8914.286
$ws.Cells.Item(88, 9).Value = 1500
$ws.Cells.Item(88, 10).Value = 10150
$ws.Cells.Item(88, 11).Value = 1500
$ws.Cells.Item(88, 12).Value = 10150
$ws.Cells.Item(88, 13).Value = -1094
$ws.Cells.Item(88, 14).Value = -10962

# ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(89, 8).Value = 6684.933
$ws.Cells.Item(89, 10).Value = 7232.5835
$ws.Cells.Item(89, 12).Value = 36162.9175
$ws.Cells.Item(89, 14).Value = -47394.9175

# ALC row 91
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(91, 8).Value = 8914.286
$ws.Cells.Item(91, 9).Value = 1500
$ws.Cells.Item(91, 10).Value = 10150
$ws.Cells.Item(91, 11).Value = 1500
$ws.Cells.Item(91, 12).Value = 10150
$ws.Cells.Item(91, 13).Value = -96
$ws.Cells.Item(91, 14).Value = -12958

# ALC row 126
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(126, 8).Value = 23000
$ws.Cells.Item(126, 10).Value = 23000
$ws.Cells.Item(126, 12).Value = 23000
$ws.Cells.Item(126, 14).Value = -32880

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(141, 8).Value = 1204.2084
$ws.Cells.Item(141, 9).Value = 1169.6086
$ws.Cells.Item(141, 10).Value = 2000
$ws.Cells.Item(141, 11).Value = 3508.8258
$ws.Cells.Item(141, 12).Value = 6000
$ws.Cells.Item(141, 13).Value = 1671.1742
$ws.Cells.Item(141, 14).Value = -16360

# ARM row 41
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(41, 8).Value = 2651.2
$ws.Cells.Item(41, 9).Value = 2651.2
$ws.Cells.Item(41, 11).Value = 2651.2
$ws.Cells.Item(41, 13).Value = -2237.2

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 5111.107
$ws.Cells.Item(86, 9).Value = 3772.9092
$ws.Cells.Item(86, 10).Value = 10017.833
$ws.Cells.Item(86, 11).Value = 3772.9092
$ws.Cells.Item(86, 12).Value = 10017.833
$ws.Cells.Item(86, 13).Value = -2649.9092
$ws.Cells.Item(86, 14).Value = -12263.833

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 5111.107
$ws.Cells.Item(89, 9).Value = 3772.9092
$ws.Cells.Item(89, 10).Value = 10017.833
$ws.Cells.Item(89, 11).Value = 18864.546
$ws.Cells.Item(89, 12).Value = 50089.165
$ws.Cells.Item(89, 13).Value = -13248.546
$ws.Cells.Item(89, 14).Value = -61321.165

# BSM row 112
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(112, 8).Value = 39965
$ws.Cells.Item(112, 10).Value = 39965
$ws.Cells.Item(112, 12).Value = 39965
$ws.Cells.Item(112, 14).Value = -42919

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 2317931.5
$ws.Cells.Item(62, 9).Value = 4276281.5
$ws.Cells.Item(62, 10).Value = 3517.6365
$ws.Cells.Item(62, 11).Value = 4276281.5
$ws.Cells.Item(62, 12).Value = 3517.6365
$ws.Cells.Item(62, 13).Value = -4275657.5
$ws.Cells.Item(62, 14).Value = -4765.636500000001

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(65, 8).Value = 2317931.5
$ws.Cells.Item(65, 9).Value = 4276281.5
$ws.Cells.Item(65, 10).Value = 3517.6365
$ws.Cells.Item(65, 11).Value = 21381407.5
$ws.Cells.Item(65, 12).Value = 17588.1825
$ws.Cells.Item(65, 13).Value = -21378287.5
$ws.Cells.Item(65, 14).Value = -23828.1825

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 1830.4286
$ws.Cells.Item(132, 9).Value = 1109.0834
$ws.Cells.Item(132, 10).Value = 3128.85
$ws.Cells.Item(132, 11).Value = 3327.2502
$ws.Cells.Item(132, 12).Value = 9386.549999999999
$ws.Cells.Item(132, 13).Value = -797.2501999999999
$ws.Cells.Item(132, 14).Value = -14446.55

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 1923.2413
$ws.Cells.Item(134, 9).Value = 1139.3334
$ws.Cells.Item(134, 11).Value = 3418.0002
$ws.Cells.Item(134, 13).Value = -883.0001999999999

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 464.73334
$ws.Cells.Item(5, 9).Value = 272
$ws.Cells.Item(5, 11).Value = 816
$ws.Cells.Item(5, 13).Value = -704

# CUL row 50
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(50, 8).Value = 145
$ws.Cells.Item(50, 9).Value = 115.55556
$ws.Cells.Item(50, 11).Value = 346.66668
$ws.Cells.Item(50, 13).Value = 134.33332

# CUL row 53
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(53, 8).Value = 145
$ws.Cells.Item(53, 9).Value = 115.55556
$ws.Cells.Item(53, 11).Value = 346.66668
$ws.Cells.Item(53, 13).Value = 134.33332

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 538.4835
$ws.Cells.Item(113, 9).Value = 531.0164
$ws.Cells.Item(113, 10).Value = 553.6667
$ws.Cells.Item(113, 11).Value = 1593.0492
$ws.Cells.Item(113, 12).Value = 1661.0001
$ws.Cells.Item(113, 13).Value = 576.9508000000001
$ws.Cells.Item(113, 14).Value = -6001.0001

# CUL row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(129, 8).Value = 2128.55
$ws.Cells.Item(129, 9).Value = 683.3333
$ws.Cells.Item(129, 10).Value = 3311
$ws.Cells.Item(129, 11).Value = 2049.9999
$ws.Cells.Item(129, 12).Value = 9933
$ws.Cells.Item(129, 13).Value = 2950.0001
$ws.Cells.Item(129, 14).Value = -19933

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(135, 8).Value = 464.73334
$ws.Cells.Item(135, 9).Value = 272
$ws.Cells.Item(135, 11).Value = 2448
$ws.Cells.Item(135, 13).Value = 87

# CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(140, 8).Value = 1204
$ws.Cells.Item(140, 9).Value = 889.0526
$ws.Cells.Item(140, 10).Value = 2700
$ws.Cells.Item(140, 11).Value = 2667.1578
$ws.Cells.Item(140, 12).Value = 8100
$ws.Cells.Item(140, 13).Value = 2512.8422
$ws.Cells.Item(140, 14).Value = -18460

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3533.122
$ws.Cells.Item(80, 9).Value = 3942.4688
$ws.Cells.Item(80, 10).Value = 2077.6667
$ws.Cells.Item(80, 11).Value = 3942.4688
$ws.Cells.Item(80, 12).Value = 2077.6667
$ws.Cells.Item(80, 13).Value = -2944.4688
$ws.Cells.Item(80, 14).Value = -4073.6667

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83, 8).Value = 3533.122
$ws.Cells.Item(83, 9).Value = 3942.4688
$ws.Cells.Item(83, 10).Value = 2077.6667
$ws.Cells.Item(83, 11).Value = 19712.344
$ws.Cells.Item(83, 12).Value = 10388.3335
$ws.Cells.Item(83, 13).Value = -14720.344
$ws.Cells.Item(83, 14).Value = -20372.3335

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 3031.456
$ws.Cells.Item(132, 9).Value = 2815.923
$ws.Cells.Item(132, 10).Value = 3498.4443
$ws.Cells.Item(132, 11).Value = 8447.769
$ws.Cells.Item(132, 12).Value = 10495.3329
$ws.Cells.Item(132, 13).Value = -5917.769
$ws.Cells.Item(132, 14).Value = -15555.3329

# LTW row 2
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 8000
$ws.Cells.Item(2, 10).Value = 8000
$ws.Cells.Item(2, 12).Value = 8000
$ws.Cells.Item(2, 14).Value = -8224

# LTW row 34
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(34, 8).Value = 7625
$ws.Cells.Item(34, 9).Value = 7625
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 11).Value = 7625
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 13).Value = -7453
$ws.Cells.Item(34, 14).ClearContents()

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 1203.2963
$ws.Cells.Item(82, 10).Value = 1167.1052
$ws.Cells.Item(82, 12).Value = 1167.1052
$ws.Cells.Item(82, 14).Value = -1889.1052

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(85, 8).Value = 1203.2963
$ws.Cells.Item(85, 10).Value = 1167.1052
$ws.Cells.Item(85, 12).Value = 1167.1052
$ws.Cells.Item(85, 14).Value = -3663.1052

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 2124.4783
$ws.Cells.Item(93, 9).Value = 2110.9333
$ws.Cells.Item(93, 11).Value = 2110.9333
$ws.Cells.Item(93, 13).Value = -862.9333000000001

# LTW row 104
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(104, 8).Value = 14515.333
$ws.Cells.Item(104, 10).Value = 14515.333
$ws.Cells.Item(104, 12).Value = 14515.333
$ws.Cells.Item(104, 14).Value = -21503.333

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 2194.7368
$ws.Cells.Item(81, 9).Value = 977.7778
$ws.Cells.Item(81, 10).Value = 3290
$ws.Cells.Item(81, 11).Value = 1955.5556
$ws.Cells.Item(81, 12).Value = 6580
$ws.Cells.Item(81, 13).Value = -894.5555999999999
$ws.Cells.Item(81, 14).Value = -8702

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(84, 8).Value = 2194.7368
$ws.Cells.Item(84, 9).Value = 977.7778
$ws.Cells.Item(84, 10).Value = 3290
$ws.Cells.Item(84, 11).Value = 9777.778
$ws.Cells.Item(84, 12).Value = 32900
$ws.Cells.Item(84, 13).Value = -4473.778
$ws.Cells.Item(84, 14).Value = -43508

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 43013156
$ws.Cells.Item(136, 9).Value = 62501932
$ws.Cells.Item(136, 10).Value = 22225128
$ws.Cells.Item(136, 11).Value = 187505796
$ws.Cells.Item(136, 12).Value = 66675384
$ws.Cells.Item(136, 13).Value = -187503246
$ws.Cells.Item(136, 14).Value = -66680484
